# JS rich text editor
# Append a new bullet line to the "Selber geschriebenes JavaScript" allowed-technologies
# cell (C8 on Tabelle1), and grow the row height to fit the extra wrapped line.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$cell = $ws.Cells.Item(8, 3)
$existingText = $cell.Value2
# Leading apostrophe keeps this a literal-text entry (cell already uses quote-prefix
# formatting since the text starts with "-") instead of it being reinterpreted.
$cell.Formula = "'" + $existingText + "`n- JS Rich Text Editor (z.B. Quill)"

$ws.Rows.Item(8).RowHeight = 75
